$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update specific C/D cell values that changed
$ws.Range("C10").Value = 601
$ws.Range("D10").Value = 601

$ws.Range("C18").Value = 531
$ws.Range("D18").Value = 531

$ws.Range("C21").Value = 608
$ws.Range("D21").Value = 608

$ws.Range("C34").Value = 218
$ws.Range("D34").Value = 274

$ws.Range("C39").Value = 24
$ws.Range("D39").Value = 402

$ws.Range("C51").Value = 2
$ws.Range("D51").Value = 88

$ws.Range("C53").Value = 308
$ws.Range("D53").Value = 308

$ws.Range("C60").Value = 959
$ws.Range("D60").Value = 959

$ws.Range("C61").Value = 526
$ws.Range("D61").Value = 526

$ws.Range("C76").Value = 274
$ws.Range("D76").Value = 274

$ws.Range("C78").Value = 402
$ws.Range("D78").Value = 402

$ws.Range("C79").Value = 424
$ws.Range("D79").Value = 424

# Delete column H entirely (header "pendientes" + all data)
$ws.Columns.Item(8).Delete()
